$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Add the new user row (row 52)
$ws.Cells.Item(52, 1).Value = "Linking_AutoUser"
$ws.Cells.Item(52, 2).Value = "Password1"
$ws.Cells.Item(52, 5).Value = "Default user for Linking tests"
$ws.Cells.Item(52, 6).Value = "N"
$ws.Cells.Item(52, 7).Value = "linking.autouser@mailinator.com"

# Apply the same border formatting used by the rest of the table
$rng = $ws.Range("A52:G52")
$rng.Borders.Color = 0

# Update the view: scroll so row 22 is at top and select E29
$ws.Activate()
[void]$ws.Range("E29").Select()
$excel.ActiveWindow.ScrollRow = 22
